# Generate Report for Archive
#
# The localization status report is being regenerated: the file that was
# previously handed off is now back "In Translation", and the Status
# columns (which are autosized to their contents) shrink to fit the new,
# shorter label.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Update the status text everywhere it appears -------------------------
# Overview sheet keeps one Status column per locale (zh-cn -> E, de-de -> F)
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"

# Each locale detail sheet has its own Status column (column C)
$wsZhCn.Range("C2").Value = "In Translation"
$wsDeDe.Range("C2").Value = "In Translation"

# --- Re-fit the Status columns to the (now shorter) text -------------------
# "In Translation" is shorter than "Ready for handoff", so the autosized
# Status columns narrow accordingly.
$wsOverview.Columns.Item(5).ColumnWidth = 12.45   # Overview!E (zh-cn status)
$wsOverview.Columns.Item(6).ColumnWidth = 12.45   # Overview!F (de-de status)
$wsZhCn.Columns.Item(3).ColumnWidth = 12.45        # zh-cn!C (Status)
$wsDeDe.Columns.Item(3).ColumnWidth = 12.45        # de-de!C (Status)
